# Applies "Add data for 2021-10-05" changes to carjacking-arrests-by-month-yoy-latest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet / update title to reflect new "through" date
$ws.Name = "Through 2021-09-27"

# Row 6 (April): 2020 arrest_made/no_arrest_made/arrest_rate columns (Q,R,S)
$ws.Range("Q6").Value = 4
$ws.Range("R6").Value = 60
$ws.Range("S6").Value = 0.0625

# Row 11 (September): update label text and 2015-2021 stats
$ws.Range("A11").Value = "September (through 09-27)"

$ws.Range("C11").Value = 27
$ws.Range("D11").Value = 0.0357

$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 37
$ws.Range("G11").Value = 0.0976

$ws.Range("I11").Value = 64
$ws.Range("J11").Value = 0.0725

$ws.Range("O11").Value = 59
$ws.Range("P11").Value = 0.0923

$ws.Range("Q11").Value = 4
$ws.Range("R11").Value = 99
$ws.Range("S11").Value = 0.0388

$ws.Range("U11").Value = 163
$ws.Range("V11").Value = 0.0121

# Row 12 (Total): update aggregated totals
$ws.Range("C12").Value = 192
$ws.Range("D12").Value = 0.1351

$ws.Range("E12").Value = 45
$ws.Range("F12").Value = 377
$ws.Range("G12").Value = 0.1066

$ws.Range("I12").Value = 570
$ws.Range("J12").Value = 0.0806

$ws.Range("O12").Value = 372
$ws.Range("P12").Value = 0.1014

$ws.Range("Q12").Value = 52
$ws.Range("R12").Value = 835
$ws.Range("S12").Value = 0.0586

$ws.Range("U12").Value = 1159
$ws.Range("V12").Value = 0.0615
